$wb = $excel.ActiveWorkbook

# Rename the existing sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Estimates"

# Add the new "Effort" sheet right after "Estimates".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Effort"

# Header row. Columns are populated in the same order shared strings were
# introduced upstream (B, C, E, F, G, then D) so the sharedStrings table
# ends up in the same sequence as the authored workbook.
$ws2.Range("A1").Value = "year"
$ws2.Range("B1").Value = "nvessels"
$ws2.Range("C1").Value = "mesh_in_avg"
$ws2.Range("E1").Value = "n_set_obs"
$ws2.Range("F1").Value = "n_days_obs"
$ws2.Range("G1").Value = "obs_perc"
$ws2.Range("D1").Value = "n_days_tot"
$ws2.Range("A1:G1").Font.Bold = $true

# Data row.
$ws2.Range("A2").Value = 2012
$ws2.Range("B2").Value = 50
$ws2.Range("C2").Value = 7.3
$ws2.Range("D2").Value = 1360
$ws2.Range("E2").Value = 250
$ws2.Range("F2").Value = 75
$ws2.Range("G2").Value = 5.5

# Column widths (approximate best-fit widths from the source workbook).
$ws2.Columns.Item(1).ColumnWidth = 4.328125
$ws2.Columns.Item(2).ColumnWidth = 7.1640625
$ws2.Columns.Item(3).ColumnWidth = 11.1640625
$ws2.Columns.Item(4).ColumnWidth = 9.328125
$ws2.Columns.Item(5).ColumnWidth = 8.6640625
$ws2.Columns.Item(6).ColumnWidth = 9.828125
$ws2.Columns.Item(7).ColumnWidth = 7.6640625

# Select the header row (matches the selection saved in the source file)
# and make "Effort" the active/visible tab.
$ws2.Rows.Item(1).Select() | Out-Null
$ws2.Activate() | Out-Null
